$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New "PERIOD TO EXPIRE" (column H) and "LAST UPDATE" (column I) values for
# rows 3-15, reflecting progress recorded as of 04-Nov-2025 (one day later
# than the previous update, so PERIOD TO EXPIRE drops by 1 for each row).
$updates = @(
    @{ Row = 3;  H = 386 },
    @{ Row = 4;  H = 288 },
    @{ Row = 5;  H = 297 },
    @{ Row = 6;  H = 694 },
    @{ Row = 7;  H = 714 },
    @{ Row = 8;  H = 290 },
    @{ Row = 9;  H = 297 },
    @{ Row = 10; H = 84 },
    @{ Row = 11; H = -98 },
    @{ Row = 12; H = 126 },
    @{ Row = 13; H = 128 },
    @{ Row = 14; H = 140 },
    @{ Row = 15; H = 612 }
)

foreach ($u in $updates) {
    # Numeric column: plain value assignment is safe.
    $ws.Cells.Item($u.Row, 8).Value = $u.H

    # Date-like text column: assigning the literal string via .Value would
    # be auto-recognized by Excel as a date and converted to a date serial
    # (changing the cell's type/number format). Instead, write it as a
    # text formula first and then paste-special just the resulting value,
    # which keeps it a plain text value and preserves the cell's existing
    # style/format untouched.
    $iCell = $ws.Cells.Item($u.Row, 9)
    $iCell.Formula = "=""04-Nov-2025"""
    $iCell.Copy()
    $iCell.PasteSpecial(-4163)
}

$excel.CutCopyMode = $false
